# Reorder the "Recorded By" (column G) value on each attendance row:
# move the leading "System" (or "admin@admin.com") token to the back,
# putting the real user identity first. Only rows whose G value is an
# exact match for one of the known before-states are touched; anything
# else (single-entry cells, already-reordered cells, etc.) is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "System, system, backup@backdoor.com" = "backup@backdoor.com, System, system"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
